$wb = $excel.ActiveWorkbook

# Add a new worksheet ("Sheet1") after the last existing sheet ("BOARD STATUS"),
# so it becomes the third / newly active tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)

# --- Populate the new sheet -------------------------------------------------
# Header row (row 1)
$ws.Range("C1").Value = "Versions"
$ws.Range("D1").Value = "Bug "
$ws.Range("B1").Value = "test"
$ws.Range("B2").Value = "CS current"
$ws.Range("D2").Value = "Test Mode entry"
$ws.Range("E1").Value = "STAUTS"
$ws.Range("E2").Value = "CLEAR"
$ws.Range("E3").Value = "REEEVIEEW"
$ws.Range("E4").Value = "NOT CHECKED"
$ws.Range("A7").Value = "BUG "
$ws.Range("B7").Value = "DEVICE"

$ws.Range("A1").Value = "Device"
$ws.Range("A2").Value = 4832
$ws.Range("A3").Value = 4832
$ws.Range("A4").Value = 4832
$ws.Range("A5").Value = 4832
$ws.Range("B3").Value = "CS current"
$ws.Range("B4").Value = "CS current"
$ws.Range("B5").Value = "CS current"
$ws.Range("C2").Value = "A0"
$ws.Range("C3").Value = "A0"
$ws.Range("C4").Value = "A0"
$ws.Range("C5").Value = "A0"
$ws.Range("D3").Value = "Test Mode entry"
$ws.Range("D4").Value = "Test Mode entry"
$ws.Range("D5").Value = "Test Mode entry"

# Bold the header row
$ws.Range("A1:E1").Font.Bold = $true

# Autofit the used columns to their content
$ws.Range("A1:F7").EntireColumn.AutoFit() | Out-Null

# Leave the same cell selected/active as in the source file
$ws.Range("E8").Select() | Out-Null
